$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H1:H3").Merge()
Write-Output "done"
